$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The workbook used to be a manually-built "Birthday Wisher" sheet with a
# hyperlinked Email column, a bordered "Status" header, and a bunch of
# pre-formatted but empty rows (4-22) left over from the original template.
# The update (commit: "whatsapp msg sender updated, now send using pywhatkit
# module") adds a "Mobile" column (needed by pywhatkit, which dials a phone
# number instead of opening an email client) and cleans the sheet down to
# just the two real data rows.
# ---------------------------------------------------------------------------

# 1) Drop the mailto: hyperlinks on the Email column (now plain text).
$ws.Hyperlinks.Delete()

# 2) Remove the leftover blank template rows (4-22) that only carried the
#    date-style formatting on column B.
$ws.Rows("4:22").Delete()

# 3) Simplify the "Status" header (F1) so it matches the other plain bold
#    headers instead of the special boxed/filled style - copy the format
#    from a neighboring plain header cell.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 4) Strip the Hyperlink character style from the former link cells so they
#    read as normal text.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2:E3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 5) Add the new "Mobile" column with the two phone numbers.
$ws.Range("G1").Value = "Mobile"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G2").Value = 1234567890
$ws.Range("G3").Value = 9876543210

# Size the new column to fit its (short) contents, the way the old DOB
# column was sized to fit its long formatted timestamps.
$ws.Columns("G:G").ColumnWidth = 10.1666666666667   # -> stored width 11

# 6) Put the selection where the author left it after entering the last
#    mobile number.
$ws.Range("B3").Select() | Out-Null
